# Y4_B2526_Excuses.xlsx - "Upload Y4_B2526_Excuses.xlsx via attendance app"
#
# The attendance app re-uploaded the sheet with a refreshed batch of rows:
#   - 7 of the 8 data rows kept the same Subject/Date/Time/Type/User but got
#     a new Student ID.
#   - Row 4's Log Time flipped from a literal text "10:30:00" to a real
#     Excel time-of-day serial (0.4375 == 10:30 AM), while row 6's Log Time
#     flipped the other way, from the serial value back to literal text.
#   - The last row (old row 9, Student ID 211926) is gone - the sheet now
#     only has 7 data rows (A1:F8 instead of A1:F9).
#
# Student IDs are stored as genuine text cells (not numbers) in the target
# file, even though they look numeric. Excel's normal Range.Value setter
# auto-converts a pure-digit string typed into a General-formatted cell
# into a Number - so a plain assignment would silently change the cell's
# type. To avoid that, Set-TextId below stages the new ID as text in an
# out-of-the-way scratch cell (Z1, far outside the A:F data range) with
# NumberFormat "@" applied, copies it, and pastes *values only*
# (xlPasteValues = -4163) into the destination. Values-only paste carries
# the text-ness of the value over without touching the destination cell's
# existing formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextId($cellAddr, $val) {
    $ws.Range("Z1").NumberFormat = "@"
    $ws.Range("Z1").Value = $val
    $ws.Range("Z1").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)   # xlPasteValues
    $ws.Range("Z1").Clear()
}

# Row 2: Student ID 212598 -> 212205 (Subject/Date/Time/Type/User unchanged)
Set-TextId "A2" "212205"

# Row 3: Student ID 211697 -> 211984
Set-TextId "A3" "211984"

# Row 4: Student ID 211769 -> 212266, and Log Time becomes a numeric
# 10:30:00 (0.4375) instead of literal text. Paste D2's formatting
# (formats only = xlPasteFormats = -4122) onto D4 so it picks up the same
# time-formatted style already used by D2, rather than minting a new one.
Set-TextId "A4" "212266"
$ws.Range("D4").Value = 0.4375
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)            # xlPasteFormats

# Row 5: Student ID 211959 -> 211757
Set-TextId "A5" "211757"

# Row 6: Student ID 212333 -> 211737, and Log Time reverts to literal text
# "10:30:00" (previously the numeric 0.4375). Borrow C6's formatting
# (same row/style family, no time number format) for D6.
Set-TextId "A6" "211737"
$ws.Range("D6").Value = "10:30:00"
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)            # xlPasteFormats

# Row 7: Student ID 211973 -> 211263
Set-TextId "A7" "211263"

# Row 8: Student ID 211687 -> 211299
Set-TextId "A8" "211299"

# Former row 9 (Student ID 211926, general surgery, 14/10/2025, 10:30:00,
# Excuse, System) is no longer present - delete it, which also shifts the
# sheet's used range/dimension down from A1:F9 to A1:F8.
$ws.Rows.Item(9).Delete()
